# Renumber the ACTB padlock/splint probe pairs (rows 4-11, column B) so the
# numbering is contiguous (0,1,2,3,4 instead of 0,2,4,6,8). This is part of
# adding overlap-checking to the probe design output: probes that overlap
# got folded out, so the surviving probes are renumbered sequentially.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value  = "ACTB_padlock_1"
$ws.Range("B5").Value  = "ACTB_splint_1"
$ws.Range("B6").Value  = "ACTB_padlock_2"
$ws.Range("B7").Value  = "ACTB_splint_2"
$ws.Range("B8").Value  = "ACTB_padlock_3"
$ws.Range("B9").Value  = "ACTB_splint_3"
$ws.Range("B10").Value = "ACTB_padlock_4"
$ws.Range("B11").Value = "ACTB_splint_4"

# Reset the view to A1 (the sheet had been left scrolled/selected at L16).
$ws.Range("A1").Select() | Out-Null

# Force a full recalculation on next load, mirroring the workbook being
# saved by an Excel build that marks the calc chain dirty after this edit.
$excel.CalculateFullRebuild() | Out-Null
